# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that comma-separated recorder names are listed in reverse order (e.g.
# "user@example.com, System" becomes "System, user@example.com").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 157

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $orig = $cell.Value2

    if ($orig -eq $null) {
        continue
    }

    $parts = $orig.Split(",")
    $count = $parts.Length

    if ($count -le 1) {
        continue
    }

    $reversed = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i].Trim()
    }

    $newVal = $reversed -join ", "
    $cell.Value = $newVal
}
